$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel COM auto-converts the assigned
# string into a numeric Value (these columns are text in the workbook).
$textCells = @("D5", "D6", "D12", "D20", "D21", "D22", "D24", "D26", "D28", "D31", "D33", "D34", "D39", "D40", "D43", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume(1h) figures scraped this run.
$ws.Range("D2").Value = "62.298.59"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.998.39"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "588.36"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "145.69"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "2.996.58"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  +3.94%  "
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "3.496.42"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "62.259.05"
$ws.Range("D19").Value = "2.999.81"
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("D20").Value = "455.63"
$ws.Range("E20").Value = "  -4.57%  "
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "81.89"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("E25").Value = "  -9.24%  "
$ws.Range("D26").Value = "12.16"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "9.72"
$ws.Range("E28").Value = "  -8.13%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "6.92"
$ws.Range("E31").Value = "  -5.62%  "
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("D33").Value = "27.57"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("D39").Value = "9.18"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "50.21"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +6.70%  "
$ws.Range("E42").Value = "  -11.63%  "
$ws.Range("D43").Value = "392.84"
$ws.Range("E43").Value = "  -9.39%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.730.37"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.266"
$ws.Range("E46").Value = "  -7.42%  "
$ws.Range("D47").Value = "37.18"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D48").Value = "129.59"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "2.18"
$ws.Range("E51").Value = "  -0.71%  "
